$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New record appended as row 10 (id 112058967 / Hygrocybe conica sighting).
$ws.Range("A10").Value = 112058967
$ws.Range("B10").Value = 86149
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 4379
$ws.Range("F10").Value = "Toppvaxskivling"
$ws.Range("G10").Value = "Hygrocybe conica"
$ws.Range("H10").Value = "(Schaeff.) P.Kumm."

$ws.Range("P10").Value = "Söder Hässleby (Söder Hässleby), Nrk"
$ws.Range("Q10").Value = 513646.2891263207
$ws.Range("R10").Value = 6578480.707588105
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = "Örebro"
$ws.Range("U10").Value = "Örebro"
$ws.Range("V10").Value = "Närke"
$ws.Range("W10").Value = "Axberg"

# Date-ish text columns: force text formatting first so Excel doesn't
# auto-convert the "yyyy-mm-dd" strings into date serial numbers, then
# drop the style back to Normal so no extra number-format survives.
$dateRange = $ws.Range("Y10:AB10")
$dateRange.NumberFormat = "@"
$ws.Range("Y10").Value = "2023-09-12"
$ws.Range("Z10").Value = "00:00"
$ws.Range("AA10").Value = "2023-09-12"
$ws.Range("AB10").Value = "00:00"
$dateRange.Style = "Normal"

$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false

$ws.Range("AW10").Value = "Erik Göthlin"
$ws.Range("AX10").Value = "Erik Göthlin"
